$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted for the "Especial" quality
# Durazno (September Snow) lot that used to sit on row 167. Insert a fresh
# row at 168 so the existing row 168 (and everything after it) shifts down
# to 169, keeping row 167 in place for now.
$ws.Rows.Item(168).Insert()

# Row 167 now becomes a brand-new record (different variety/quality/origin).
$ws.Range("D167").Value = 44509
$ws.Range("K167").Value = "Florida King"
$ws.Range("L167").Value = "Tercera"
$ws.Range("M167").Value = 250
$ws.Range("N167").Value = 18000
$ws.Range("O167").Value = 18000
$ws.Range("P167").Value = 18000
$ws.Range("R167").Value = "Provincia de Limarí"
$ws.Range("S167").Value = 1286

# The freshly inserted row 168 is populated with the "Especial" record that
# previously lived on row 167, preserving the market/region/product context.
$ws.Range("A168").Value = 4
$ws.Range("B168").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C168").Value = "Los Lagos"
$ws.Range("D168").Value = 44273
$ws.Range("E168").Value = 10
$ws.Range("F168").Value = "Fruta"
$ws.Range("G168").Value = 100103
$ws.Range("H168").Value = "Frutos de hueso (carozo)"
$ws.Range("I168").Value = 100103004
$ws.Range("J168").Value = "Durazno"
$ws.Range("K168").Value = "September Snow"
$ws.Range("L168").Value = "Especial"
$ws.Range("M168").Value = 150
$ws.Range("N168").Value = 20000
$ws.Range("O168").Value = 20000
$ws.Range("P168").Value = 20000
$ws.Range("Q168").Value = "$/caja 14 kilos empedrada"
$ws.Range("R168").Value = "Región de O'Higgins"
$ws.Range("S168").Value = 1429
$ws.Range("T168").Value = 14
